# MTRN4230 Asst 2 Requirements Tracking - "Updating Requirement Sheet [Draft]"
#
# Applies the authors' edits to the "Status" worksheet:
#  - bump the "Last updated" date
#  - record the updater ("Ken") and group letter ("L") in the header block
#  - fill in "Responsible person" (column E) for a batch of requirement rows
#  - widen columns D and F to fit the new responsible-person text
#  - move the active cell selection to C4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status")

# --- header block (rows 2-3) ---------------------------------------------
$ws.Range("B2").Value = 20170824
$ws.Range("B3").Value = "Ken"
$ws.Range("D3").Value = "L"

# --- column widths (D: Date last tested was too narrow; F: Tested by) ----
# Target stored widths are 17.140625 / 40.28515625 "character" units; the
# ColumnWidth setter here only resolves to sixths of a character, so these
# inputs land on the closest representable width (17.1667 / 40.3333).
$ws.Columns.Item(4).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 39.5

# --- Responsible person (column E) ----------------------------------------
$responsible = @{
    8  = "Ravi,Emtiazul,Jay"
    9  = "Ravi,Emtiazul,Jay"
    10 = "Ravi,Emtiazul,Jay"
    11 = "Ravi,Emtiazul,Jay"
    12 = "Ravi,Emtiazul,Jay"
    13 = "Ravi,Emtiazul,Jay"
    14 = "Undetermined"
    15 = "Jay"
    18 = "Daniel, Herman, Ken"
    19 = "Daniel, Herman, Ken"
    22 = "Daniel, Herman, Ken"
    23 = "Herman"
    24 = "Daniel, Herman, Ken"
    27 = "Undetermined"
    28 = "Undetermined"
    29 = "Undetermined"
    30 = "Undetermined"
    31 = "Undetermined"
    32 = "Undetermined"
    33 = "Undetermined"
    36 = "Everyone"
    37 = "Everyone"
    40 = "Ken"
    41 = "Everyone"
    42 = "Everyone"
}

foreach ($row in $responsible.Keys) {
    $ws.Cells.Item($row, 5).Value = $responsible[$row]
}

# --- selection --------------------------------------------------------------
$ws.Range("C4").Select()
